$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates (F1 renamed, G1:K1 added with header style copied from F1) ---
$ws.Cells.Item(1,6).Value = "frequency"
$headerMap = [ordered]@{
  7  = "frequency_occurrence"
  8  = "frequency_occurrence_probab"
  9  = "max_probab"
  10 = "max_probab_percentage"
  11 = "recommended_level"
}
foreach ($col in $headerMap.Keys) {
  $ws.Range("F1").Copy($ws.Cells.Item(1, $col))
  $ws.Cells.Item(1, $col).Value = $headerMap[$col]
}

# --- Pre-materialize every cell in the new G:K range (rows 2-60) so that rows which stay
#     blank (no frequency data) still get an explicit (empty) cell, matching the source rows
#     that already contained empty placeholder cells in columns A-F. ---
$ws.Range("G2:I60").NumberFormat = "General"
$ws.Range("G2:I60").Style = "Normal"
$ws.Range("K2:K60").NumberFormat = "General"
$ws.Range("K2:K60").Style = "Normal"

# --- Force column J (max_probab_percentage) to be stored as text, even though its values look numeric ---
$ws.Range("J2:J60").NumberFormat = "@"

# --- Data rows ---
# row 2
$ws.Cells.Item(2,6).Value = 9
$ws.Cells.Item(2,7).Value = '{"L3":9}'
$ws.Cells.Item(2,8).Value = '{"L3":1.0}'
$ws.Cells.Item(2,9).Value = 1
$ws.Cells.Item(2,10).Value = '100.00'
$ws.Cells.Item(2,11).Value = 'L3'
# row 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = '{"L3":1}'
$ws.Cells.Item(3,8).Value = '{"L3":1.0}'
$ws.Cells.Item(3,9).Value = 1
$ws.Cells.Item(3,10).Value = '100.00'
$ws.Cells.Item(3,11).Value = 'L3'
# row 4
$ws.Cells.Item(4,6).Value = 7
$ws.Cells.Item(4,7).Value = '{"L3":7}'
$ws.Cells.Item(4,8).Value = '{"L3":1.0}'
$ws.Cells.Item(4,9).Value = 1
$ws.Cells.Item(4,10).Value = '100.00'
$ws.Cells.Item(4,11).Value = 'L3'
# row 5
$ws.Cells.Item(5,6).Value = 2
$ws.Cells.Item(5,7).Value = '{"L3":2}'
$ws.Cells.Item(5,8).Value = '{"L3":1.0}'
$ws.Cells.Item(5,9).Value = 1
$ws.Cells.Item(5,10).Value = '100.00'
$ws.Cells.Item(5,11).Value = 'L3'
# row 6
$ws.Cells.Item(6,6).Value = 7
$ws.Cells.Item(6,7).Value = '{"L3":7}'
$ws.Cells.Item(6,8).Value = '{"L3":1.0}'
$ws.Cells.Item(6,9).Value = 1
$ws.Cells.Item(6,10).Value = '100.00'
$ws.Cells.Item(6,11).Value = 'L3'
# row 7
$ws.Cells.Item(7,6).Value = 9
$ws.Cells.Item(7,7).Value = '{"L1":6,"L2":3}'
$ws.Cells.Item(7,8).Value = '{"L1":0.6666666667,"L2":0.3333333333}'
$ws.Cells.Item(7,9).Value = 0.6666666666666666
$ws.Cells.Item(7,10).Value = '66.67'
$ws.Cells.Item(7,11).Value = 'L2'
# row 8
$ws.Cells.Item(8,6).Value = 2
$ws.Cells.Item(8,7).Value = '{"L3":2}'
$ws.Cells.Item(8,8).Value = '{"L3":1.0}'
$ws.Cells.Item(8,9).Value = 1
$ws.Cells.Item(8,10).Value = '100.00'
$ws.Cells.Item(8,11).Value = 'L3'
# row 9
$ws.Cells.Item(9,6).Value = 6
$ws.Cells.Item(9,7).Value = '{"L3":6}'
$ws.Cells.Item(9,8).Value = '{"L3":1.0}'
$ws.Cells.Item(9,9).Value = 1
$ws.Cells.Item(9,10).Value = '100.00'
$ws.Cells.Item(9,11).Value = 'L3'
# row 10
$ws.Cells.Item(10,6).Value = 6
$ws.Cells.Item(10,7).Value = '{"L3":6}'
$ws.Cells.Item(10,8).Value = '{"L3":1.0}'
$ws.Cells.Item(10,9).Value = 1
$ws.Cells.Item(10,10).Value = '100.00'
$ws.Cells.Item(10,11).Value = 'L3'
# row 11
$ws.Cells.Item(11,6).Value = 7
$ws.Cells.Item(11,7).Value = '{"L3":6,"L2":1}'
$ws.Cells.Item(11,8).Value = '{"L3":0.8571428571,"L2":0.1428571429}'
$ws.Cells.Item(11,9).Value = 0.8571428571428571
$ws.Cells.Item(11,10).Value = '85.71'
$ws.Cells.Item(11,11).Value = 'L3'
# row 12
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = '{"L2":1}'
$ws.Cells.Item(12,8).Value = '{"L2":1.0}'
$ws.Cells.Item(12,9).Value = 1
$ws.Cells.Item(12,10).Value = '100.00'
$ws.Cells.Item(12,11).Value = 'L2'
# row 13
$ws.Cells.Item(13,6).Value = 9
$ws.Cells.Item(13,7).Value = '{"L3":6,"L2":3}'
$ws.Cells.Item(13,8).Value = '{"L3":0.6666666667,"L2":0.3333333333}'
$ws.Cells.Item(13,9).Value = 0.6666666666666666
$ws.Cells.Item(13,10).Value = '66.67'
$ws.Cells.Item(13,11).Value = 'L3'
# row 14
$ws.Cells.Item(14,6).Value = 9
$ws.Cells.Item(14,7).Value = '{"L3":6,"L2":3}'
$ws.Cells.Item(14,8).Value = '{"L3":0.6666666667,"L2":0.3333333333}'
$ws.Cells.Item(14,9).Value = 0.6666666666666666
$ws.Cells.Item(14,10).Value = '66.67'
$ws.Cells.Item(14,11).Value = 'L3'
# row 15
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = '{"L1":1}'
$ws.Cells.Item(15,8).Value = '{"L1":1.0}'
$ws.Cells.Item(15,9).Value = 1
$ws.Cells.Item(15,10).Value = '100.00'
$ws.Cells.Item(15,11).Value = 'L1'
# row 16
$ws.Cells.Item(16,6).Value = 7
$ws.Cells.Item(16,7).Value = '{"L3":7}'
$ws.Cells.Item(16,8).Value = '{"L3":1.0}'
$ws.Cells.Item(16,9).Value = 1
$ws.Cells.Item(16,10).Value = '100.00'
$ws.Cells.Item(16,11).Value = 'L3'
# row 17
$ws.Cells.Item(17,6).Value = 8
$ws.Cells.Item(17,7).Value = '{"L2":4,"L3":2,"L1":2}'
$ws.Cells.Item(17,8).Value = '{"L2":0.5,"L3":0.25,"L1":0.25}'
$ws.Cells.Item(17,9).Value = 0.5
$ws.Cells.Item(17,10).Value = '50.00'
$ws.Cells.Item(17,11).Value = 'L3'
# row 18
$ws.Cells.Item(18,6).Value = 8
$ws.Cells.Item(18,7).Value = '{"L2":4,"L3":2,"L1":2}'
$ws.Cells.Item(18,8).Value = '{"L2":0.5,"L3":0.25,"L1":0.25}'
$ws.Cells.Item(18,9).Value = 0.5
$ws.Cells.Item(18,10).Value = '50.00'
$ws.Cells.Item(18,11).Value = 'L3'
# row 19
$ws.Cells.Item(19,6).Value = 7
$ws.Cells.Item(19,7).Value = '{"L2":6,"L1":1}'
$ws.Cells.Item(19,8).Value = '{"L2":0.8571428571,"L1":0.1428571429}'
$ws.Cells.Item(19,9).Value = 0.8571428571428571
$ws.Cells.Item(19,10).Value = '85.71'
$ws.Cells.Item(19,11).Value = 'L2'
# row 20
$ws.Cells.Item(20,6).Value = 9
$ws.Cells.Item(20,7).Value = '{"L3":9}'
$ws.Cells.Item(20,8).Value = '{"L3":1.0}'
$ws.Cells.Item(20,9).Value = 1
$ws.Cells.Item(20,10).Value = '100.00'
$ws.Cells.Item(20,11).Value = 'L3'
# row 21
$ws.Cells.Item(21,6).Value = 1
$ws.Cells.Item(21,7).Value = '{"L3":1}'
$ws.Cells.Item(21,8).Value = '{"L3":1.0}'
$ws.Cells.Item(21,9).Value = 1
$ws.Cells.Item(21,10).Value = '100.00'
$ws.Cells.Item(21,11).Value = 'L3'
# row 22
$ws.Cells.Item(22,6).Value = 7
$ws.Cells.Item(22,7).Value = '{"L2":6,"L3":1}'
$ws.Cells.Item(22,8).Value = '{"L2":0.8571428571,"L3":0.1428571429}'
$ws.Cells.Item(22,9).Value = 0.8571428571428571
$ws.Cells.Item(22,10).Value = '85.71'
$ws.Cells.Item(22,11).Value = 'L3'
# row 23
# row 24
$ws.Cells.Item(24,6).Value = 9
$ws.Cells.Item(24,7).Value = '{"L3":9}'
$ws.Cells.Item(24,8).Value = '{"L3":1.0}'
$ws.Cells.Item(24,9).Value = 1
$ws.Cells.Item(24,10).Value = '100.00'
$ws.Cells.Item(24,11).Value = 'L3'
# row 25
$ws.Cells.Item(25,6).Value = 9
$ws.Cells.Item(25,7).Value = '{"L2":9}'
$ws.Cells.Item(25,8).Value = '{"L2":1.0}'
$ws.Cells.Item(25,9).Value = 1
$ws.Cells.Item(25,10).Value = '100.00'
$ws.Cells.Item(25,11).Value = 'L2'
# row 26
$ws.Cells.Item(26,6).Value = 9
$ws.Cells.Item(26,7).Value = '{"L2":6,"L3":2,"L1":1}'
$ws.Cells.Item(26,8).Value = '{"L2":0.6666666667,"L3":0.2222222222,"L1":0.1111111111}'
$ws.Cells.Item(26,9).Value = 0.6666666666666666
$ws.Cells.Item(26,10).Value = '66.67'
$ws.Cells.Item(26,11).Value = 'L3'
# row 27
$ws.Cells.Item(27,6).Value = 10
$ws.Cells.Item(27,7).Value = '{"L2":7,"L1":3}'
$ws.Cells.Item(27,8).Value = '{"L2":0.7,"L1":0.3}'
$ws.Cells.Item(27,9).Value = 0.7
$ws.Cells.Item(27,10).Value = '70.00'
$ws.Cells.Item(27,11).Value = 'L2'
# row 28
$ws.Cells.Item(28,6).Value = 9
$ws.Cells.Item(28,7).Value = '{"L2":7,"L1":2}'
$ws.Cells.Item(28,8).Value = '{"L2":0.7777777778,"L1":0.2222222222}'
$ws.Cells.Item(28,9).Value = 0.7777777777777778
$ws.Cells.Item(28,10).Value = '77.78'
$ws.Cells.Item(28,11).Value = 'L2'
# row 29
# row 30
$ws.Cells.Item(30,6).Value = 15
$ws.Cells.Item(30,7).Value = '{"L2":7,"L3":7,"L1":1}'
$ws.Cells.Item(30,8).Value = '{"L2":0.4666666667,"L3":0.4666666667,"L1":0.0666666667}'
$ws.Cells.Item(30,9).Value = 0.4666666666666667
$ws.Cells.Item(30,10).Value = '46.67'
$ws.Cells.Item(30,11).Value = 'L3'
# row 31
$ws.Cells.Item(31,6).Value = 10
$ws.Cells.Item(31,7).Value = '{"L3":10}'
$ws.Cells.Item(31,8).Value = '{"L3":1.0}'
$ws.Cells.Item(31,9).Value = 1
$ws.Cells.Item(31,10).Value = '100.00'
$ws.Cells.Item(31,11).Value = 'L3'
# row 32
# row 33
# row 34
$ws.Cells.Item(34,6).Value = 11
$ws.Cells.Item(34,7).Value = '{"L3":11}'
$ws.Cells.Item(34,8).Value = '{"L3":1.0}'
$ws.Cells.Item(34,9).Value = 1
$ws.Cells.Item(34,10).Value = '100.00'
$ws.Cells.Item(34,11).Value = 'L3'
# row 35
$ws.Cells.Item(35,6).Value = 11
$ws.Cells.Item(35,7).Value = '{"L3":11}'
$ws.Cells.Item(35,8).Value = '{"L3":1.0}'
$ws.Cells.Item(35,9).Value = 1
$ws.Cells.Item(35,10).Value = '100.00'
$ws.Cells.Item(35,11).Value = 'L3'
# row 36
$ws.Cells.Item(36,6).Value = 9
$ws.Cells.Item(36,7).Value = '{"L2":5,"L1":4}'
$ws.Cells.Item(36,8).Value = '{"L2":0.5555555556,"L1":0.4444444444}'
$ws.Cells.Item(36,9).Value = 0.5555555555555556
$ws.Cells.Item(36,10).Value = '55.56'
$ws.Cells.Item(36,11).Value = 'L2'
# row 37
$ws.Cells.Item(37,6).Value = 9
$ws.Cells.Item(37,7).Value = '{"L2":8,"L3":1}'
$ws.Cells.Item(37,8).Value = '{"L2":0.8888888889,"L3":0.1111111111}'
$ws.Cells.Item(37,9).Value = 0.8888888888888888
$ws.Cells.Item(37,10).Value = '88.89'
$ws.Cells.Item(37,11).Value = 'L3'
# row 38
$ws.Cells.Item(38,6).Value = 10
$ws.Cells.Item(38,7).Value = '{"L2":8,"L3":2}'
$ws.Cells.Item(38,8).Value = '{"L2":0.8,"L3":0.2}'
$ws.Cells.Item(38,9).Value = 0.8
$ws.Cells.Item(38,10).Value = '80.00'
$ws.Cells.Item(38,11).Value = 'L3'
# row 39
$ws.Cells.Item(39,6).Value = 10
$ws.Cells.Item(39,7).Value = '{"L1":6,"L2":4}'
$ws.Cells.Item(39,8).Value = '{"L1":0.6,"L2":0.4}'
$ws.Cells.Item(39,9).Value = 0.6
$ws.Cells.Item(39,10).Value = '60.00'
$ws.Cells.Item(39,11).Value = 'L2'
# row 40
# row 41
# row 42
# row 43
# row 44
$ws.Cells.Item(44,6).Value = 1
$ws.Cells.Item(44,7).Value = '{"L3":1}'
$ws.Cells.Item(44,8).Value = '{"L3":1.0}'
$ws.Cells.Item(44,9).Value = 1
$ws.Cells.Item(44,10).Value = '100.00'
$ws.Cells.Item(44,11).Value = 'L3'
# row 45
$ws.Cells.Item(45,6).Value = 9
$ws.Cells.Item(45,7).Value = '{"L1":8,"L2":1}'
$ws.Cells.Item(45,8).Value = '{"L1":0.8888888889,"L2":0.1111111111}'
$ws.Cells.Item(45,9).Value = 0.8888888888888888
$ws.Cells.Item(45,10).Value = '88.89'
$ws.Cells.Item(45,11).Value = 'L2'
# row 46
$ws.Cells.Item(46,6).Value = 9
$ws.Cells.Item(46,7).Value = '{"L1":9}'
$ws.Cells.Item(46,8).Value = '{"L1":1.0}'
$ws.Cells.Item(46,9).Value = 1
$ws.Cells.Item(46,10).Value = '100.00'
$ws.Cells.Item(46,11).Value = 'L1'
# row 47
$ws.Cells.Item(47,6).Value = 9
$ws.Cells.Item(47,7).Value = '{"L1":9}'
$ws.Cells.Item(47,8).Value = '{"L1":1.0}'
$ws.Cells.Item(47,9).Value = 1
$ws.Cells.Item(47,10).Value = '100.00'
$ws.Cells.Item(47,11).Value = 'L1'
# row 48
# row 49
# row 50
# row 51
# row 52
# row 53
# row 54
# row 55
# row 56
$ws.Cells.Item(56,6).Value = 9
$ws.Cells.Item(56,7).Value = '{"L2":6,"L3":2,"L1":1}'
$ws.Cells.Item(56,8).Value = '{"L2":0.6666666667,"L3":0.2222222222,"L1":0.1111111111}'
$ws.Cells.Item(56,9).Value = 0.6666666666666666
$ws.Cells.Item(56,10).Value = '66.67'
$ws.Cells.Item(56,11).Value = 'L3'
# row 57
$ws.Cells.Item(57,6).Value = 7
$ws.Cells.Item(57,7).Value = '{"L3":7}'
$ws.Cells.Item(57,8).Value = '{"L3":1.0}'
$ws.Cells.Item(57,9).Value = 1
$ws.Cells.Item(57,10).Value = '100.00'
$ws.Cells.Item(57,11).Value = 'L3'
# row 58
$ws.Cells.Item(58,6).Value = 10
$ws.Cells.Item(58,7).Value = '{"L1":10}'
$ws.Cells.Item(58,8).Value = '{"L1":1.0}'
$ws.Cells.Item(58,9).Value = 1
$ws.Cells.Item(58,10).Value = '100.00'
$ws.Cells.Item(58,11).Value = 'L1'
# row 59
$ws.Cells.Item(59,6).Value = 9
$ws.Cells.Item(59,7).Value = '{"L2":6,"L3":3}'
$ws.Cells.Item(59,8).Value = '{"L2":0.6666666667,"L3":0.3333333333}'
$ws.Cells.Item(59,9).Value = 0.6666666666666666
$ws.Cells.Item(59,10).Value = '66.67'
$ws.Cells.Item(59,11).Value = 'L3'
# row 60
$ws.Cells.Item(60,6).Value = 9
$ws.Cells.Item(60,7).Value = '{"L1":8,"L2":1}'
$ws.Cells.Item(60,8).Value = '{"L1":0.8888888889,"L2":0.1111111111}'
$ws.Cells.Item(60,9).Value = 0.8888888888888888
$ws.Cells.Item(60,10).Value = '88.89'
$ws.Cells.Item(60,11).Value = 'L2'

# Reset style for column J cells we forced to text, to drop the temporary number-format style
$ws.Range("J2:J60").Style = "Normal"

Write-Output "edit complete"
